$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# text so Excel does not auto-convert the string into a numeric value.
$textCells = @('D5', 'D6', 'D13', 'D14', 'D18', 'D19', 'D22', 'D24', 'D27', 'D28', 'D30', 'D33', 'D37', 'D40', 'D43', 'D45', 'D47', 'D50')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.667.56'
$ws.Range('E2').Value = '  -4.30%  '
$ws.Range('D3').Value = '3.340.45'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('D5').Value = '573.66'
$ws.Range('E5').Value = '  -3.31%  '
$ws.Range('D6').Value = '180.89'
$ws.Range('E6').Value = '  -5.34%  '
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('E11').Value = '  -3.42%  '
$ws.Range('D12').Value = '3.918.40'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').Value = '0.136'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '27.03'
$ws.Range('E14').Value = '  -5.58%  '
$ws.Range('D15').Value = '66.760.72'
$ws.Range('E15').Value = '  -4.16%  '
$ws.Range('E16').Value = '  -2.77%  '
$ws.Range('D17').Value = '3.323.80'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = '437.43'
$ws.Range('E18').Value = '  -3.56%  '
$ws.Range('D19').Value = '5.70'
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('D22').Value = '73.56'
$ws.Range('E22').Value = '  -3.09%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Value = '0.518'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('E25').Value = '  -3.87%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '9.04'
$ws.Range('E27').Value = '  -4.96%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('D30').Value = '22.83'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  -5.77%  '
$ws.Range('D33').Value = '6.77'
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('E36').Value = '  -6.08%  '
$ws.Range('D37').Value = '27.31'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -8.36%  '
$ws.Range('D39').Value = '2.835.94'
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').Value = '0.798'
$ws.Range('E40').Value = '  -1.90%  '
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('E42').Value = '  -5.52%  '
$ws.Range('D43').Value = '40.22'
$ws.Range('E43').Value = '  -2.15%  '
$ws.Range('E44').Value = '  -3.08%  '
$ws.Range('D45').Value = '24.42'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('E46').Value = '  -7.07%  '
$ws.Range('D47').Value = '321.48'
$ws.Range('E47').Value = '  -5.30%  '
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '0.977'
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('E51').Value = '  -2.60%  '

# Restore the default (Normal) style on the cells we forced to text so
# no stray style index is left attached to them.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
